# Kr stripped vs unstripped comparison
$wb = $excel.ActiveWorkbook

# 1. Rename the first worksheet
$wsConservative = $wb.Worksheets.Item("2_no_PS_bunch_splitting_conserv")
$wsConservative.Name = "1_baseline_conservative"

# 2. Reset the scroll position on that sheet back to the top-left (A1),
#    while keeping the existing selection (F27).
$wsConservative.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$wsConservative.Range("F27").Select()

# 3. Update the selected cell on the optimistic sheet from B26 to D45.
$wsOptimistic = $wb.Worksheets.Item("2_no_PS_bunch_splitting_optimis")
$wsOptimistic.Activate()
$wsOptimistic.Range("D45").Select()
